# BINGO_cc.xlsx update
# - Added and tweaked a few phrases in the "list" sheet (the master phrase
#   list that the "mon"/"tue"/.../"mon1"/"grid" boards pull from).
# - Fixed a typo ("jouney" -> "journey").
# - Removed the stray duplicate "Can you reach out to  ____?" phrase and
#   replaced it with the corrected "Can you reach out to  _____?" wording.
# NOTE: the author's absolute local file path recorded by Excel/AutoRecover
# (x15ac:absPath in xl/workbook.xml) is editing-machine metadata that Excel
# stamps on save; it is not reachable through the Workbook/Worksheet COM
# object model exposed here, so it is not set by this script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# Order matters: the underlying shared-string table appends newly
# introduced text at the end in the order cells are (re)written, so these
# assignments are ordered to match the canonical edit sequence.

# Fix typo: "jouney" -> "journey" (row 42, bottom of the current list)
$ws.Range("A42").Value = "[on a / It's been a] journey"

# Tweak existing phrases
$ws.Range("A5").Value = "[We] can't see your screen"
$ws.Range("A11").Value = "[Let's] take this offline"

# Append new phrases
$ws.Range("A43").Value = "You're breaking up again"
$ws.Range("A44").Value = "I'm having computer problems"
$ws.Range("A12").Value = "Dog barking"
$ws.Range("A45").Value = "on the same page"

# Replace the old duplicate phrase with the corrected one (extra underscore)
$ws.Range("A27").Value = "Can you reach out to  _____?"

# Update the sheet's view: drop the scrolled-down topLeftCell and move the
# active selection from A43 to B32.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B32").Select()
